$d = $word.ActiveDocument

# Part 1: merge the " " and "affyring)" runs into a single run " affyring)"
$d.Content.Find.Execute(" affyring)", $true, $false, $false, $false, $false, $true, 1, $false, " affyring)", 2) | Out-Null

# Part 2: replace the final (bookmark-only) paragraph with the new Iteration 1 section.
# InsertXML replaces the contents of the target range with the supplied WordprocessingML fragment.
$p = $d.Paragraphs.Last
$r = $p.Range
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Overskrift2"/></w:pPr><w:r><w:t>Iteration 1:</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Need</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> to have:</w:t></w:r></w:p><w:p><w:sdt><w:sdtPr><w:id w:val="927930864"/><w14:checkbox><w14:checked w14:val="0"/><w14:checkedState w14:val="2612" w14:font="MS Gothic"/><w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/></w14:checkbox></w:sdtPr><w:sdtContent><w:r><w:rPr><w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/></w:rPr><w:t>&#9744;</w:t></w:r></w:sdtContent></w:sdt><w:r><w:t xml:space="preserve"> main menu</w:t></w:r></w:p><w:p><w:sdt><w:sdtPr><w:id w:val="1477024845"/><w14:checkbox><w14:checked w14:val="0"/><w14:checkedState w14:val="2612" w14:font="MS Gothic"/><w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/></w14:checkbox></w:sdtPr><w:sdtContent><w:r><w:rPr><w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/></w:rPr><w:t>&#9744;</w:t></w:r></w:sdtContent></w:sdt><w:r><w:t xml:space="preserve"></w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gamescreen</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:sdt><w:sdtPr><w:id w:val="765960026"/><w14:checkbox><w14:checked w14:val="0"/><w14:checkedState w14:val="2612" w14:font="MS Gothic"/><w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/></w14:checkbox></w:sdtPr><w:sdtContent><w:r><w:rPr><w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/></w:rPr><w:t>&#9744;</w:t></w:r></w:sdtContent></w:sdt><w:r><w:t xml:space="preserve"> pausefunktion</w:t></w:r></w:p><w:p><w:sdt><w:sdtPr><w:id w:val="-1079432246"/><w14:checkbox><w14:checked w14:val="0"/><w14:checkedState w14:val="2612" w14:font="MS Gothic"/><w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/></w14:checkbox></w:sdtPr><w:sdtContent><w:r><w:rPr><w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/></w:rPr><w:t>&#9744;</w:t></w:r></w:sdtContent></w:sdt><w:r><w:t xml:space="preserve"> finde modeller/baggrunde</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Nice to have:</w:t></w:r></w:p><w:p><w:sdt><w:sdtPr><w:id w:val="1061449900"/><w14:checkbox><w14:checked w14:val="0"/><w14:checkedState w14:val="2612" w14:font="MS Gothic"/><w14:uncheckedState w14:val="2610" w14:font="MS Gothic"/></w14:checkbox></w:sdtPr><w:sdtContent><w:r><w:rPr><w:rFonts w:ascii="MS Gothic" w:eastAsia="MS Gothic" w:hAnsi="MS Gothic" w:hint="eastAsia"/></w:rPr><w:t>&#9744;</w:t></w:r></w:sdtContent></w:sdt><w:r><w:t xml:space="preserve"> implementere </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>spaceship</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p/><w:p/><w:p/><w:p/><w:p/>
'@
$r.InsertXML($xml)
